$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.578.09'
$ws.Range("E2").Value = '  +0.23%  '

$ws.Range("D3").Value = '1.644.81'
$ws.Range("E3").Value = '  -0.92%  '

$ws.Range("E4").Value = '  +0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.55%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.532'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.35%  '

$ws.Range("E7").Value = '  +0.20%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.29'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.60%  '

$ws.Range("E9").Value = '  -2.74%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0610'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.37%  '

$ws.Range("E11").Value = '  +1.40%  '

$ws.Range("D12").Value = '1.877.54'
$ws.Range("E12").Value = '  -0.88%  '

$ws.Range("D13").Value = '1.648.61'
$ws.Range("E13").Value = '  -0.65%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.584'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.12%  '

$ws.Range("E15").Value = '  -2.52%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.37'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.30%  '

$ws.Range("D17").Value = '27.539.58'
$ws.Range("E17").Value = '  +0.08%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.48'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.57%  '

$ws.Range("D19").Value = '0.0₃0722'
$ws.Range("E19").Value = '  -0.94%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.58'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.32%  '

$ws.Range("E21").Value = '  +0.07%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.32'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.93%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.70'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.14%  '

$ws.Range("E24").Value = '  -2.29%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.79'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.95%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.99'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.15%  '

$ws.Range("E27").Value = '  +1.24%  '

$ws.Range("E28").Value = '  +0.31%  '

$ws.Range("E29").Value = '  -4.39%  '

$ws.Range("E30").Value = '  -0.76%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0487'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.45%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.30'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.56%  '

$ws.Range("E33").Value = '  +2.89%  '

$ws.Range("D34").Value = '1.426.20'
$ws.Range("E34").Value = '  -2.12%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.60'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.60%  '

$ws.Range("E36").Value = '  -0.53%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.570'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.45%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.883'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.45%  '

$ws.Range("E39").Value = '  -2.20%  '

$ws.Range("E40").Value = '  -2.59%  '

$ws.Range("E41").Value = '  +0.25%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.818'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.72%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.46'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.26%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.23'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.50%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '65.08'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.51%  '

$ws.Range("D46").Value = '1.786.72'
$ws.Range("E46").Value = '  -0.74%  '

$ws.Range("E47").Value = '  -2.23%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.16'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.50%  '

$ws.Range("E49").Value = '  +1.13%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0998'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.26%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.76'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.92%  '
